$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-12 14:16:37"
$wsZh.Range("E5").Value = "2016-03-12 14:16:37"
$wsZh.Range("H4").Value = "2016-03-12 14:16:55"
$wsZh.Range("H5").Value = "2016-03-12 14:16:55"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-12 14:16:40"
$wsDe.Range("E5").Value = "2016-03-12 14:16:40"
$wsDe.Range("H4").Value = "2016-03-12 14:17:01"
$wsDe.Range("H5").Value = "2016-03-12 14:17:01"
